# Update "Periodo Mora" values in column E (rows 16-19) to reflect the
# reordering of older estados de cuenta (swap row16<->row19, row17<->row18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value2 = "2203"
$ws.Range("E17").Value2 = "2202"
$ws.Range("E18").Value2 = "2201"
$ws.Range("E19").Value2 = "2112"

# Update "Valor Mora" values in column F (rows 16 and 19 swap)
$ws.Range("F16").Value2 = 10902
$ws.Range("F19").Value2 = 36341
